$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. TC02_RegisterAccount: fill in the new rows 2-6 (Email / Country code /
#    phone number / Plate number / Password) below the existing row 1.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("TC02_RegisterAccount")

# Row 2 - Email (hyperlink-styled value, added after the quote-prefixed
# cells below so the "Hyperlink" cell style lands at cellXfs index 3).
$ws2.Range("A2").Value = "Email "

# Row 3 - Country code
$ws2.Range("A3").Value = "Country code"
$ws2.Range("B3").Value = "Australia"

# Row 4 - phone number (quote-prefixed numeric-looking text)
$ws2.Range("A4").Value = "phone number"

# Row 5 - Plate number
$ws2.Range("A5").Value = "Plate number"
$ws2.Range("B5").Value = "PT00000001"

# Row 6 - Password (quote-prefixed)
$ws2.Range("A6").Value = "Password"
$ws2.Range("B6").Value = "'123456789oO"

# Now write the e-mail value + hyperlink (creates the "Hyperlink" style),
# then the quote-prefixed phone number.
$ws2.Range("B2").Value = "user317@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:user317@gmail.com")
$ws2.Range("B4").Value = "'491258667"

# Row heights to match the authored layout.
$ws2.Rows.Item(2).RowHeight = 48.75
$ws2.Rows.Item(3).RowHeight = 48
$ws2.Rows.Item(4).RowHeight = 48.75
$ws2.Rows.Item(5).RowHeight = 51.75
$ws2.Rows.Item(6).RowHeight = 59.25

# Selection moves to B4, and this sheet is no longer the active tab.
$ws2.Range("B4").Select()

# ---------------------------------------------------------------------------
# 2. Insert a new worksheet "TC03_UpdatePassword" right after
#    TC02_RegisterAccount (and before Resources), becoming the active tab.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "TC03_UpdatePassword"

$ws3.Columns.Item(1).ColumnWidth = 17.42578125
$ws3.Columns.Item(2).ColumnWidth = 45.28515625

# Row 2 - Old password
$ws3.Range("A2").Value = "Old password"

# Row 3 - New password
$ws3.Range("A3").Value = "New password"
$ws3.Range("B3").Value = "123456789oO"

# Row 1 - Email (hyperlink-styled)
$ws3.Range("A1").Value = "Email"

$ws3.Range("B2").Value = "123456789aA"

# E-mail value + hyperlink for row 1.
$ws3.Range("B1").Value = "user300@gmail.com"
$ws3.Hyperlinks.Add($ws3.Range("B1"), "mailto:user300@gmail.com")

$ws3.Rows.Item(1).RowHeight = 49.5
$ws3.Rows.Item(2).RowHeight = 48
$ws3.Rows.Item(3).RowHeight = 45.75

$ws3.Range("B3").Select()
$ws3.Activate()

Write-Output "done"
